$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the current column C (BeginDateTime), shifting
# BeginDateTime/EndDateTime one column to the right.
$ws.Columns("C").Insert()

# Header for the newly inserted column.
$ws.Range("C1").Value = "IsDefaultNationality"

# Values for the new column per row. A leading apostrophe forces literal
# text entry so "True"/"False" are stored as text instead of being
# auto-coerced to boolean TRUE/FALSE; ClearFormats then drops the
# quote-prefix formatting marker that text entry leaves behind, restoring
# the default cell style.
$ws.Range("C2").Value = "'True"
$ws.Range("C2").ClearFormats()

$ws.Range("C3").Value = "'False"
$ws.Range("C3").ClearFormats()
